$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update column F (想去人数) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 63
$ws1.Range("F3").Value = 362
$ws1.Range("F4").Value = 434
$ws1.Range("F6").Value = 13540
$ws1.Range("F7").Value = 64
$ws1.Range("F8").Value = 63
$ws1.Range("F9").Value = 5486
$ws1.Range("F12").Value = 28
$ws1.Range("F14").Value = 1215
$ws1.Range("F15").Value = 60
$ws1.Range("F16").Value = 158
$ws1.Range("F17").Value = 727
$ws1.Range("F18").Value = 2893
$ws1.Range("F19").Value = 7840
$ws1.Range("F20").Value = 1176
$ws1.Range("F22").Value = 233

# Sheet "演出" (sheet2): update column F value
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 32

# Sheet "全部类型" (sheet4): update column F (想去人数) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 63
$ws4.Range("F3").Value = 362
$ws4.Range("F4").Value = 32
$ws4.Range("F5").Value = 434
$ws4.Range("F7").Value = 13540
$ws4.Range("F8").Value = 64
$ws4.Range("F9").Value = 63
$ws4.Range("F10").Value = 5486
$ws4.Range("F13").Value = 28
$ws4.Range("F15").Value = 1215
$ws4.Range("F16").Value = 60
$ws4.Range("F17").Value = 158
$ws4.Range("F18").Value = 727
$ws4.Range("F19").Value = 2893
$ws4.Range("F21").Value = 7840
$ws4.Range("F22").Value = 1176
$ws4.Range("F24").Value = 233
